# Adjusted time trends for retirement and education processes.
# The "Parameters" sheet had a row for MAX_AGE_TO_ENTER_EDUCATION (row 9)
# removed, and the row above it (row 8) renamed from
# MAX_AGE_TO_LEAVE_CONTINUOUS_EDUCATION to
# MAX_AGE_TO_STAY_IN_CONTINUOUS_EDUCATION (keeping its value of 29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Rename the KEY in row 8 (value in B8 stays at 29).
$ws.Range("A8").Value = "MAX_AGE_TO_STAY_IN_CONTINUOUS_EDUCATION"

# Delete the whole row 9 (MAX_AGE_TO_ENTER_EDUCATION), shifting rows 10-41 up.
$ws.Rows.Item(9).Delete()

# Move the active selection to A25, matching the post-edit workbook state.
$ws.Activate()
$ws.Range("A25").Select()
